# Auto-committed on 2023/02/10 週五 17:24:00.48
# Adds a new "訂正處理" worksheet documenting the AcHcode/EntAc/TitaHCode
# correction-handling workflow, and makes it the active sheet.

$wb = $excel.ActiveWorkbook

# --- Add the new worksheet at the end, named 訂正處理 -----------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "訂正處理"

$ws = $wb.Worksheets.Item("訂正處理")

# --- Column widths -----------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 28.44140625
$ws.Columns.Item(2).ColumnWidth = 12.109375
$ws.Columns.Item(3).ColumnWidth = 21.5546875
$ws.Columns.Item(4).ColumnWidth = 28.77734375
$ws.Columns.Item(5).ColumnWidth = 22.33203125
$ws.Columns.Item(6).ColumnWidth = 13.109375
$ws.Columns.Item(7).ColumnWidth = 24.6640625
$ws.Columns.Item(8).ColumnWidth = 20.88671875

# --- Header notes -------------------------------------------------------------
$ws.Range("A1").Value = "AcHcode 帳務訂正記號 0-正常 1-刪除帳務2-沖正帳務(沖正、要入帳) 3-沖正帳務(訂正、不入帳)"
$ws.Range("A2").Value = "EntAc 入總帳記號 0:未入帳 1:已入帳"
$ws.Range("A3").Value = "TitaHCode 訂正別 0:正常 1:訂正 2:被訂正 3:沖正 4:被沖正5:兩段式已放行訂正"

# --- Table header row (row 4) -------------------------------------------------
$ws.Range("C4").Value = "pattern"
$ws.Range("D4").Value = "AcHcode "
$ws.Range("E4").Value = "分錄"
$ws.Range("F4").Value = "EntAc"
$ws.Range("G4").Value = "TitaHCode"
$ws.Range("H4").Value = "更新會計相關檔"

# --- Row 5: 1.經辦登錄 ---------------------------------------------------------
$ws.Range("A5").Value = "1.經辦登錄"
$ws.Range("C5").Value = "1.經辦登錄"
$ws.Range("D5").Value = "0-正常"
$ws.Range("E5").Value = "產生正向帳務"
$ws.Range("F5").Value = " 0:未入帳"
$ws.Range("G5").Value = "0:正常"
$ws.Range("H5").Value = "v"

# --- Row 7-8: 2.經辦訂正 -------------------------------------------------------
$ws.Range("A7").Value = "2.經辦訂正"
$ws.Range("C7").Value = "2.經辦訂正"
$ws.Range("D7").Value = "3-沖正帳務(訂正、不入帳)"
$ws.Range("E7").Value = "更新正向帳務"
$ws.Range("F7").Value = " 0:未入帳"
$ws.Range("G7").Value = "2:被訂正"
$ws.Range("H7").Value = "v"

$ws.Range("E8").Value = "產生反向帳務"
$ws.Range("F8").Value = " 0:未入帳"
$ws.Range("G8").Value = "1:訂正"
$ws.Range("H8").Value = "x"

# --- Row 10: 3.經辦修正 / 3.修正刪除 --------------------------------------------
$ws.Range("A10").Value = "3.經辦修正 "
$ws.Range("B10").Value = " 1).訂正 "
$ws.Range("C10").Value = "3.修正刪除"
$ws.Range("D10").Value = "1-刪除帳務"
$ws.Range("E10").Value = "刪除正向帳務"
$ws.Range("H10").Value = "v"

# --- Row 12: 2).登錄 ----------------------------------------------------------
$ws.Range("B12").Value = " 2).登錄"
$ws.Range("C12").Value = "1.經辦登錄"
$ws.Range("D12").Value = "0-正常"
$ws.Range("E12").Value = "產生正向帳務"
$ws.Range("F12").Value = " 0:未入帳"
$ws.Range("G12").Value = "0:正常"
$ws.Range("H12").Value = "v"

# --- Row 14-15: 4.主管放行 -----------------------------------------------------
$ws.Range("A14").Value = "4.主管放行"
$ws.Range("C14").Value = "4.主管放行"
$ws.Range("D14").Value = "0-正常"
$ws.Range("E14").Value = "更新正向帳務"
$ws.Range("F14").Value = "1:已入帳"
$ws.Range("G14").Value = "0:正常"
$ws.Range("H14").Value = "v"

$ws.Range("A15").Value = [string][char]96

# --- Row 16: 5.放行訂正 ---------------------------------------------------------
$ws.Range("A16").Value = "5.放行訂正"
$ws.Range("C16").Value = "5.放行訂正"
$ws.Range("D16").Value = "1-刪除帳務"
$ws.Range("E16").Value = "更新正向帳務"
$ws.Range("F16").Value = "0:未入帳"
$ws.Range("G16").Value = "5:兩段式已放行訂正"
$ws.Range("H16").Value = "v"

# --- Row 18-19: 6.已放行訂正、經辦修正 -------------------------------------------
$ws.Range("A18").Value = "6.已放行訂正、經辦修正 "
$ws.Range("B18").Value = " 1).訂正 "
$ws.Range("C18").Value = "2.經辦訂正"
$ws.Range("D18").Value = "3-沖正帳務(訂正、不入帳)"
$ws.Range("E18").Value = "更新正向帳務"
$ws.Range("F18").Value = " 0:未入帳"
$ws.Range("G18").Value = "2:被訂正"
$ws.Range("H18").Value = "v"

$ws.Range("E19").Value = "產生反向帳務"
$ws.Range("F19").Value = " 0:未入帳"
$ws.Range("G19").Value = "1:訂正"
$ws.Range("H19").Value = "x"

# --- Row 21: 2).登錄 -----------------------------------------------------------
$ws.Range("B21").Value = " 2).登錄"
$ws.Range("C21").Value = "1.經辦登錄"
$ws.Range("D21").Value = "0-正常"
$ws.Range("E21").Value = "產生正向帳務"
$ws.Range("F21").Value = " 0:未入帳"
$ws.Range("G21").Value = "0:正常"
$ws.Range("H21").Value = "v"

# --- Row 23-24: 7.已放行訂正、經辦訂正 -------------------------------------------
$ws.Range("A23").Value = "7.已放行訂正、經辦訂正"
$ws.Range("C23").Value = "2.經辦訂正"
$ws.Range("D23").Value = "3-沖正帳務(訂正、不入帳)"
$ws.Range("E23").Value = "更新正向帳務"
$ws.Range("F23").Value = " 0:未入帳"
$ws.Range("G23").Value = "2:被訂正"
$ws.Range("H23").Value = "v"

$ws.Range("E24").Value = "產生反向帳務"
$ws.Range("F24").Value = " 0:未入帳"
$ws.Range("G24").Value = "1:訂正"
$ws.Range("H24").Value = "x"

# --- Column C formatting (style idx 62 - left/center aligned, spans rows 4-27) -
$ws.Range("C4:C27").HorizontalAlignment = -4131  # xlLeft
$ws.Range("C4:C27").VerticalAlignment = -4108    # xlCenter

# --- Sheet view: scroll so row 4 is the top visible row, select H21 ------------
$ws.Activate()
$ws.Range("H21").Select()
$excel.ActiveWindow.ScrollRow = 4

$wb.Save()
